$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsAll = $wb.Worksheets.Item("全部类型")

# 展览
$wsExhibit.Cells.Item(2, 6).Value = 131
$wsExhibit.Cells.Item(3, 6).Value = 191
$wsExhibit.Cells.Item(4, 6).Value = 443
$wsExhibit.Cells.Item(5, 6).Value = 208
$wsExhibit.Cells.Item(7, 6).Value = 1237
$wsExhibit.Cells.Item(8, 6).Value = 413
$wsExhibit.Cells.Item(13, 6).Value = 429
$wsExhibit.Cells.Item(14, 6).Value = 802
$wsExhibit.Cells.Item(16, 6).Value = 736
$wsExhibit.Cells.Item(17, 6).Value = 298
$wsExhibit.Cells.Item(19, 6).Value = 1029
$wsExhibit.Cells.Item(20, 6).Value = 481
$wsExhibit.Cells.Item(21, 6).Value = 279
$wsExhibit.Cells.Item(22, 6).Value = 95
$wsExhibit.Cells.Item(23, 6).Value = 392
$wsExhibit.Cells.Item(26, 6).Value = 486
$wsExhibit.Cells.Item(27, 6).Value = 14

# 演出
$wsShow.Cells.Item(4, 6).Value = 371
$wsShow.Cells.Item(5, 6).Value = 45
$wsShow.Cells.Item(6, 6).Value = 46
$wsShow.Cells.Item(11, 6).Value = 153
$wsShow.Cells.Item(12, 6).Value = 97
$wsShow.Cells.Item(13, 6).Value = 34

# 本地生活
$wsLocal.Cells.Item(2, 6).Value = 350

# 全部类型
$wsAll.Cells.Item(2, 6).Value = 350
$wsAll.Cells.Item(4, 6).Value = 131
$wsAll.Cells.Item(5, 6).Value = 191
$wsAll.Cells.Item(6, 6).Value = 443
$wsAll.Cells.Item(7, 6).Value = 208
$wsAll.Cells.Item(9, 6).Value = 1237
$wsAll.Cells.Item(10, 6).Value = 413
$wsAll.Cells.Item(14, 6).Value = 371
$wsAll.Cells.Item(16, 6).Value = 45
$wsAll.Cells.Item(18, 6).Value = 46
$wsAll.Cells.Item(20, 6).Value = 429
$wsAll.Cells.Item(21, 6).Value = 802
$wsAll.Cells.Item(23, 6).Value = 736
$wsAll.Cells.Item(24, 6).Value = 298
$wsAll.Cells.Item(26, 6).Value = 1029
$wsAll.Cells.Item(27, 6).Value = 481
$wsAll.Cells.Item(30, 6).Value = 279
$wsAll.Cells.Item(31, 6).Value = 95
$wsAll.Cells.Item(32, 6).Value = 392
$wsAll.Cells.Item(34, 6).Value = 153
$wsAll.Cells.Item(37, 6).Value = 97
$wsAll.Cells.Item(38, 6).Value = 34
$wsAll.Cells.Item(39, 6).Value = 486
$wsAll.Cells.Item(42, 6).Value = 14
